# Refresh of the stock-quote export: fixes a bug where "Ações
# Internacionais" and "Renda Fixa" sheets were exported with zeroed
# values, and refreshes the latest quotes/quantities/timestamps for
# all sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: Acoes Nacionais ---
$ws1 = $wb.Worksheets.Item("Ações Nacionais")
$ws1.Range("B2").Value = 12.30646991729736
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 12.30646991729736
$ws1.Range("E2").Value = "01/12/2025 09:04:44"
$ws1.Range("B3").Value = 22.46999931335449
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 22.46999931335449
$ws1.Range("E3").Value = "01/12/2025 09:04:44"
$ws1.Range("B4").Value = 13.85999965667725
$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 13.85999965667725
$ws1.Range("E4").Value = "01/12/2025 09:04:44"
$ws1.Range("B5").Value = 35.43000030517578
$ws1.Range("C5").Value = 1
$ws1.Range("D5").Value = 35.43000030517578
$ws1.Range("E5").Value = "01/12/2025 09:04:45"
$ws1.Range("B6").Value = 67.40000152587891
$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = 67.40000152587891
$ws1.Range("E6").Value = "01/12/2025 09:04:45"
$ws1.Range("B7").Value = 19.43000030517578
$ws1.Range("C7").Value = 1
$ws1.Range("D7").Value = 19.43000030517578
$ws1.Range("E7").Value = "01/12/2025 09:04:45"
$ws1.Range("B8").Value = 11.5
$ws1.Range("C8").Value = 1
$ws1.Range("D8").Value = 11.5
$ws1.Range("E8").Value = "01/12/2025 09:04:45"
$ws1.Range("B9").Value = 30.63999938964844
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 30.63999938964844
$ws1.Range("E9").Value = "01/12/2025 09:04:45"
$ws1.Range("B10").Value = 41.63000106811523
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 41.63000106811523
$ws1.Range("E10").Value = "01/12/2025 09:04:45"
$ws1.Range("B11").Value = 7.360000133514404
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 7.360000133514404
$ws1.Range("E11").Value = "01/12/2025 09:04:46"

# --- Sheet: Acoes Internacionais ---
$ws2 = $wb.Worksheets.Item("Ações Internacionais")
$ws2.Range("A1").Value = "Ação:"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 5.334199905395508
$ws2.Range("D2").Value = 5.334199905395508
$ws2.Range("E2").Value = "01/12/2025 09:04:50"
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 5.334199905395508
$ws2.Range("D3").Value = 5.334199905395508
$ws2.Range("E3").Value = "01/12/2025 09:04:50"
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = 5.334199905395508
$ws2.Range("D4").Value = 5.334199905395508
$ws2.Range("E4").Value = "01/12/2025 09:04:50"
$ws2.Range("B5").Value = 1
$ws2.Range("C5").Value = 5.334199905395508
$ws2.Range("D5").Value = 5.334199905395508
$ws2.Range("E5").Value = "01/12/2025 09:04:50"

# --- Sheet: FIIs ---
$ws3 = $wb.Worksheets.Item("FIIs")
$ws3.Range("B2").Value = 151.0399932861328
$ws3.Range("C2").Value = 1
$ws3.Range("D2").Value = 151.0399932861328
$ws3.Range("E2").Value = "01/12/2025 09:04:58"
$ws3.Range("B3").Value = 87.25
$ws3.Range("C3").Value = 1
$ws3.Range("D3").Value = 87.25
$ws3.Range("E3").Value = "01/12/2025 09:04:58"
$ws3.Range("B4").Value = 9.920000076293945
$ws3.Range("C4").Value = 1
$ws3.Range("D4").Value = 9.920000076293945
$ws3.Range("E4").Value = "01/12/2025 09:04:58"
$ws3.Range("B5").Value = 83.04000091552734
$ws3.Range("C5").Value = 1
$ws3.Range("D5").Value = 83.04000091552734
$ws3.Range("E5").Value = "01/12/2025 09:04:59"
$ws3.Range("B6").Value = 109.5800018310547
$ws3.Range("C6").Value = 1
$ws3.Range("D6").Value = 109.5800018310547
$ws3.Range("E6").Value = "01/12/2025 09:04:59"
$ws3.Range("B7").Value = 79.23000335693359
$ws3.Range("C7").Value = 1
$ws3.Range("D7").Value = 79.23000335693359
$ws3.Range("E7").Value = "01/12/2025 09:04:59"
$ws3.Range("B8").Value = 79.87000274658203
$ws3.Range("C8").Value = 1
$ws3.Range("D8").Value = 79.87000274658203
$ws3.Range("E8").Value = "01/12/2025 09:04:59"
$ws3.Range("B9").Value = 88.27999877929688
$ws3.Range("C9").Value = 1
$ws3.Range("D9").Value = 88.27999877929688
$ws3.Range("E9").Value = "01/12/2025 09:04:59"
$ws3.Range("C10").Value = 1
$ws3.Range("D10").Value = 64.48999786376953
$ws3.Range("E10").Value = "01/12/2025 09:04:59"
$ws3.Range("C11").Value = 1
$ws3.Range("D11").Value = 8.380000114440918
$ws3.Range("E11").Value = "01/12/2025 09:05:00"
$ws3.Range("B12").Value = 5.25
$ws3.Range("C12").Value = 1
$ws3.Range("D12").Value = 5.25
$ws3.Range("E12").Value = "01/12/2025 09:05:00"
$ws3.Range("B13").Value = 9.630000114440918
$ws3.Range("C13").Value = 1
$ws3.Range("D13").Value = 9.630000114440918
$ws3.Range("E13").Value = "01/12/2025 09:05:00"

# --- Sheet: Renda Fixa ---
$ws4 = $wb.Worksheets.Item("Renda Fixa")

# New header columns D1 and E1 need the same bold/boxed header style as
# the existing header cells, so copy formatting from C1 (already styled).
$ws4.Range("C1").Copy()
$ws4.Range("D1:E1").PasteSpecial(-4122)

# New D2 (Valor em Reais) needs the currency number style used by B2.
$ws4.Range("B2").Copy()
$ws4.Range("D2").PasteSpecial(-4122)

# Remove the old B2/C2 values (columns no longer used on this row).
$ws4.Range("B2").Clear()
$ws4.Range("C2").Clear()

# Update header labels
$ws4.Range("A1").Value = "Ação"
$ws4.Range("B1").Value = "Valor em USD"
$ws4.Range("C1").Value = "Cotação Dolar"
$ws4.Range("D1").Value = "Valor em Reais"
$ws4.Range("E1").Value = "Atualizado em"

# Update data row
$ws4.Range("D2").Value = 100
$ws4.Range("E2").Value = "01/12/2025 09:05:04"
